$wb = $excel.ActiveWorkbook

# Rename sheets (order corresponds to sheetId 1-5)
$wb.Worksheets.Item(1).Name = "GNG_TO-1650477881132394"
$wb.Worksheets.Item(2).Name = "NB_TO-1650477883796042"
$wb.Worksheets.Item(3).Name = "RS_TO-16504778837980092"
$wb.Worksheets.Item(4).Name = "TOL_TO-16504778838440104"
$wb.Worksheets.Item(5).Name = "vSAT_TO-16504778839050398"

# Sheet 1 (GNG)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-16504778811043937.csv"
$ws1.Range("B3").Value = "GNG_stims-16504778811154273.csv"
$ws1.Range("B4").Value = "go_stims-16504778811163926.csv"
$ws1.Range("B5").Value = "GNG_stims-1650477881131392.csv"

# Sheet 2 (NB)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "ZB-match_7-1650477881339391.csv"
$ws2.Range("B3").Value = "ZB-match_0-16504778814823954.csv"
$ws2.Range("B4").Value = "TB-16504778826793938.csv"
$ws2.Range("B5").Value = "OB-16504778818443933.csv"
$ws2.Range("B6").Value = "ZB-match_3-16504778814593966.csv"
$ws2.Range("B7").Value = "TB-16504778820384257.csv"
$ws2.Range("B8").Value = "TB-16504778837690392.csv"
$ws2.Range("B9").Value = "OB-16504778816013942.csv"
$ws2.Range("B10").Value = "OB-16504778818603961.csv"

# Sheet 4 (TOL)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-1650477883812007.csv"
$ws4.Range("B3").Value = "ZM_stims-1650477883799009.csv"
$ws4.Range("B4").Value = "MM_stims-16504778838280413.csv"
$ws4.Range("B5").Value = "ZM_stims-16504778838130088.csv"
$ws4.Range("B6").Value = "MM_stims-165047788384304.csv"
$ws4.Range("B7").Value = "ZM_stims-16504778838290062.csv"

# Sheet 5 (vSAT)
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "vSAT_stims-1650477883874006.csv"
$ws5.Range("B3").Value = "vSAT_stims-16504778838910065.csv"
$ws5.Range("B4").Value = "SAT_stims-165047788384601.csv"
$ws5.Range("B5").Value = "SAT_stims-16504778838580086.csv"
